# Regen sval data to filter save games
# Update columns B, C, D, E, G for rows 2-10 on the active sheet (F is unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 0.6753301551942219;  C = 1.667794583268128;  D = 0.1575252929769615;  E = 0.496779210170732;  G = 2.997429241610044 }
    3  = @{ B = 0.6753301551942219;  C = 0.3127903958511391;  D = 3.900430680208489;   E = 0.496779210170732;  G = 5.385330441424582 }
    4  = @{ B = 1.459612070389937;   C = 1.667794583268128;  D = 26.21740644021617;   E = 8.660232485948974;  G = 38.00504557982321 }
    5  = @{ B = 0.6753301551942219;  C = 1.667794583268128;  D = 0.8054896365839992;  E = 0.496779210170732;  G = 3.645393585217082 }
    6  = @{ B = 0.04763786555579896; C = 0.04240448674262143; D = 0.8054896365839992;  E = 8.660232485948974;  G = 9.555764474831394 }
    7  = @{ B = 0.127881588408715;   C = 0.3127903958511391;  D = 0.8054896365839992;  E = 8.660232485948974;  G = 9.906394106792828 }
    8  = @{ B = 1.459612070389937;   C = 1.667794583268128;  D = 0.1575252929769615;  E = 0.496779210170732;  G = 3.781711156805759 }
    9  = @{ B = 3.230985683306322;   C = 1.667794583268128;  D = 0.8054896365839992;  E = 0.496779210170732;  G = 6.201049113329182 }
    10 = @{ B = 0.6753301551942219;  C = 1.667794583268128;  D = 0.8054896365839992;  E = 8.660232485948974;  G = 11.80884686099532 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
